$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Existing-row tweaks (rows 17, 22, 25, 27-30)
# ---------------------------------------------------------------------------

# Row 17: fill in "BUG FIXED BY" (M17) with "Hemant", bold font to match style 32
$ws.Range("M17").Value = "Hemant"
$ws.Range("M17").Font.Bold = $true

# Row 22: count corrections
$ws.Range("F22").Value = 32
$ws.Range("H22").Value = 32

# Row 25: count corrections
$ws.Range("D25").Value = 109
$ws.Range("E25").Value = 86
$ws.Range("F25").Value = 23
$ws.Range("H25").Value = 23

# Row 30: count corrections + status + formatting fixes
$ws.Range("D30").Value = 67
$ws.Range("E30").Value = 54
$ws.Range("F30").Value = 13
$ws.Range("H30").Value = 13

# K30 / Q30 had stray one-off styles (50 / 51) in the source file; bring them
# back in line with the rest of the column by copying the format used by the
# equivalent cells elsewhere in the table.
$ws.Range("K25").Copy() | Out-Null
$ws.Range("K30").PasteSpecial(-4122) | Out-Null
$ws.Range("K30").Value = 46066

$ws.Range("Q16").Copy() | Out-Null
$ws.Range("Q30").PasteSpecial(-4122) | Out-Null

$ws.Range("A16").Copy() | Out-Null
$ws.Range("N30").PasteSpecial(-4122) | Out-Null
$ws.Range("N30").Value = "Pending"

$ws.Range("A16").Copy() | Out-Null
$ws.Range("P30").PasteSpecial(-4122) | Out-Null
$ws.Range("P30").Value = "OtherInvoice_Search"

# ---------------------------------------------------------------------------
# 2) Brand-new rows 31-33
# ---------------------------------------------------------------------------

# --- Row 31 ------------------------------------------------------------
$ws.Range("A30:S30").Copy() | Out-Null
$ws.Range("A31:S31").PasteSpecial(-4122) | Out-Null

$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "Operation"
$ws.Range("C31").Value = "GST Invoice->Other Invoice>Print"
$ws.Range("D31").Value = 52
$ws.Range("E31").Value = 32
$ws.Range("F31").Value = 20
$ws.Range("G31").ClearContents() | Out-Null
$ws.Range("H31").Value = 20
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = "Akash varun"
$ws.Range("K31").Value = 46069
$ws.Range("L31").ClearContents() | Out-Null
$ws.Range("M31").ClearContents() | Out-Null
$ws.Range("N31").Value = "Pending"
$ws.Range("O31").Value = "(this bug ID is used because the same issue occurs. - 713,719)"
$ws.Range("P31").Value = "OtherInvoice_Print"
$ws.Range("Q31").Formula = '=HYPERLINK("[Icaffe_Export Sea_Operaton Test Cases.xlsx]OtherInvoice_Print!A1","Go to  OtherInvoice_Print")'
$ws.Range("R31").Value = "Sudhir kumar sinha"
$ws.Range("S31").ClearContents() | Out-Null

$ws.Range("C31").Copy() | Out-Null
$ws.Range("D27").Copy() | Out-Null
$ws.Range("C31").PasteSpecial(-4122) | Out-Null
$ws.Range("C31").Value = "GST Invoice->Other Invoice>Print"

$ws.Range("O23").Copy() | Out-Null
$ws.Range("O31").PasteSpecial(-4122) | Out-Null
$ws.Range("O31").Value = "(this bug ID is used because the same issue occurs. - 713,719)"

$ws.Rows(31).RowHeight = 45

# --- Row 32 ------------------------------------------------------------
$ws.Range("A31:S31").Copy() | Out-Null
$ws.Range("A32:S32").PasteSpecial(-4122) | Out-Null

$ws.Range("A32").Value = 31
$ws.Range("B32").Value = "Operation"
$ws.Range("C32").Value = "GST Invoice->Other Invoice>Print"
$ws.Range("D32").Value = 2
$ws.Range("E32").ClearContents() | Out-Null
$ws.Range("F32").ClearContents() | Out-Null
$ws.Range("G32").Value = 2
$ws.Range("H32").ClearContents() | Out-Null
$ws.Range("I32").ClearContents() | Out-Null
$ws.Range("J32").Value = "Akash varun"
$ws.Range("K32").Value = 46069
$ws.Range("L32").ClearContents() | Out-Null
$ws.Range("M32").ClearContents() | Out-Null
$ws.Range("N32").Value = "Pending"
$ws.Range("P32").Value = "OtherInvoice_Print"
$ws.Range("Q32").Formula = '=HYPERLINK("[Icaffe_Export Sea_Operaton Test Cases.xlsx]OtherInvoice_Print!A1","Go to  OtherInvoice_Print")'
$ws.Range("R32").Value = "Sudhir kumar sinha"
$ws.Range("S32").ClearContents() | Out-Null

$ws.Range("O22").Copy() | Out-Null
$ws.Range("O32").PasteSpecial(-4122) | Out-Null
$ws.Range("O32").ClearContents() | Out-Null

$ws.Rows(32).RowHeight = 30

# --- Row 33 ------------------------------------------------------------
$ws.Range("A30:S30").Copy() | Out-Null
$ws.Range("A33:S33").PasteSpecial(-4122) | Out-Null

$ws.Range("A33").Value = 32
$ws.Range("B33").Value = "Operation"
$ws.Range("C33").Value = "GST Invoice->Consultancy Invoice"
$ws.Range("D33").Value = 93
$ws.Range("E33").Value = 78
$ws.Range("F33").Value = 15
$ws.Range("G33").ClearContents() | Out-Null
$ws.Range("H33").Value = 15
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = "Akash varun"
$ws.Range("K33").Value = 46070
$ws.Range("L33").ClearContents() | Out-Null
$ws.Range("M33").ClearContents() | Out-Null
$ws.Range("N33").Value = "Working"
$ws.Range("O33").ClearContents() | Out-Null
$ws.Range("P33").Value = "Consultancy_Invoice"
$ws.Range("Q33").Formula = '=HYPERLINK("[Icaffe_Export Sea_Operaton Test Cases.xlsx]Consultancy_Invoice!A1","Go to  Consultancy_Invoice")'
$ws.Range("R33").Value = "Sudhir kumar sinha"
$ws.Range("S33").ClearContents() | Out-Null

$ws.Range("D27").Copy() | Out-Null
$ws.Range("C33").PasteSpecial(-4122) | Out-Null
$ws.Range("C33").Value = "GST Invoice->Consultancy Invoice"

$ws.Range("B23").Copy() | Out-Null
$ws.Range("B33").PasteSpecial(-4122) | Out-Null
$ws.Range("B33").Value = "Operation"

$ws.Rows(33).RowHeight = 30

# ---------------------------------------------------------------------------
# 3) Sheet view bookkeeping: scroll/freeze pane & selection moved with the
#    newly-added rows.
# ---------------------------------------------------------------------------
$ws.Range("H33").Select() | Out-Null
$ws.Application.ActiveWindow.ScrollRow = 24
